$wb = $excel.ActiveWorkbook

# --- Update StatOutput sheet: number_of_files -> 37, number_of_sample -> 3 ---
# Values must stay stored as shared strings (text), not numbers, and without
# picking up an extra cell style, so we force text format, assign, then
# clear the formatting that was just applied (keeps the string type, drops
# the style index).
$statOutput = $wb.Worksheets.Item("StatOutput")

$statOutput.Range("A2:B2").NumberFormat = "@"
$statOutput.Range("A2").Value = "37"
$statOutput.Range("B2").Value = "3"
$statOutput.Range("A2:B2").ClearFormats()

# --- Update StatOutput_Message sheet: Cypher query text now filters on
#     'Cocker Spaniel' instead of 'Akita' ---
$statMessage = $wb.Worksheets.Item("StatOutput_Message")
$cypher = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN['Cocker Spaniel']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"
$statMessage.Range("A18").Value = $cypher
